# Applies the 07-01-2024 scraper update to the Cambodia CPL 2023-2024 sheet:
#  - rows 38/39, 40/41 and 58/59 each had their match data (columns F, H..V)
#    swapped (the fixture order changed upstream, index/date/country/league
#    columns A-E and the home-goals column G stay put);
#  - three brand-new fixtures are appended as rows 79-81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param($rowA, $rowB, $bufferRow)

    # Use a scratch row far outside the used range as a temp buffer so we
    # don't have to hand-retype every field.
    $ws.Range("F${rowA}:V${rowA}").Copy() | Out-Null
    $ws.Range("F${bufferRow}:V${bufferRow}").PasteSpecial(-4104) | Out-Null

    $ws.Range("F${rowB}:V${rowB}").Copy() | Out-Null
    $ws.Range("F${rowA}:V${rowA}").PasteSpecial(-4104) | Out-Null

    $ws.Range("F${bufferRow}:V${bufferRow}").Copy() | Out-Null
    $ws.Range("F${rowB}:V${rowB}").PasteSpecial(-4104) | Out-Null

    $ws.Range("F${bufferRow}:V${bufferRow}").ClearContents() | Out-Null
}

$excel.CutCopyMode = 0

Swap-MatchRows 38 39 200
Swap-MatchRows 40 41 200
Swap-MatchRows 58 59 200

$excel.CutCopyMode = 0

function Add-MatchRow {
    param(
        $r,
        $indice, $pais, $torneio, $temporada, $dataPartida,
        $home, $homeGols, $away, $awayGols,
        $homeOpenOdds, $homeOpenData, $homeCloseOdds, $homeCloseData,
        $drawOpenOdds, $drawOpenData, $drawCloseOdds, $drawCloseData,
        $awayOpenOdds, $awayOpenData, $awayCloseOdds, $awayCloseData,
        $url
    )

    # Copy formatting (bold/border style on A, date style on E) from the
    # last existing data row before filling in the new values.
    $ws.Range("A78:V78").Copy() | Out-Null
    $ws.Range("A${r}:V${r}").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    $ws.Range("A$r").Value = $indice
    $ws.Range("B$r").Value = $pais
    $ws.Range("C$r").Value = $torneio
    $ws.Range("D$r").Value = $temporada
    $ws.Range("E$r").Value = $dataPartida
    $ws.Range("F$r").Value = $home
    $ws.Range("G$r").Value = $homeGols
    $ws.Range("H$r").Value = $away
    $ws.Range("I$r").Value = $awayGols
    $ws.Range("J$r").Value = $homeOpenOdds
    $ws.Range("K$r").Value = $homeOpenData
    $ws.Range("L$r").Value = $homeCloseOdds
    $ws.Range("M$r").Value = $homeCloseData
    $ws.Range("N$r").Value = $drawOpenOdds
    $ws.Range("O$r").Value = $drawOpenData
    $ws.Range("P$r").Value = $drawCloseOdds
    $ws.Range("Q$r").Value = $drawCloseData
    $ws.Range("R$r").Value = $awayOpenOdds
    $ws.Range("S$r").Value = $awayOpenData
    $ws.Range("T$r").Value = $awayCloseOdds
    $ws.Range("U$r").Value = $awayCloseData
    $ws.Range("V$r").Value = $url
}

Add-MatchRow 79 78 "cambodia" "cpl" "2023-2024" 45298.5 `
    "Svay Rieng" 2 "Prey Veng" 0 `
    1.22 "07/01/2024 01:12" 1.15 "07/01/2024 11:33" `
    5.67 "07/01/2024 01:12" 7.43 "07/01/2024 11:56" `
    8.279999999999999 "07/01/2024 01:12" 9.91 "07/01/2024 11:56" `
    "https://www.betexplorer.com/football/cambodia/cpl/svay-rieng-prey-veng/z3MfFgzC/"

Add-MatchRow 80 79 "cambodia" "cpl" "2023-2024" 45298.5 `
    "Phnom Penh Crown" 3 "Angkor Tiger" 2 `
    1.16 "07/01/2024 01:12" 1.15 "07/01/2024 11:58" `
    6.73 "07/01/2024 01:12" 7.33 "07/01/2024 11:59" `
    9.52 "07/01/2024 01:12" 9.970000000000001 "07/01/2024 11:59" `
    "https://www.betexplorer.com/football/cambodia/cpl/phnom-penh-crown-angkor-tiger/EgLbEDkI/"

Add-MatchRow 81 80 "cambodia" "cpl" "2023-2024" 45298.5 `
    "Visakha" 1 "Tiffy Army" 2 `
    1.79 "07/01/2024 02:12" 1.47 "07/01/2024 11:58" `
    3.51 "07/01/2024 02:12" 4.33 "07/01/2024 11:58" `
    3.68 "07/01/2024 02:12" 5.1 "07/01/2024 11:58" `
    "https://www.betexplorer.com/football/cambodia/cpl/visakha-tiffy-army/fi5cGZK5/"

$excel.CutCopyMode = 0
